# "Wires up level end script"
#
# The LevelCompleteTitle row ("Great!") had its German and Polish cells
# both mistakenly pointing at the English "Super!" string. Give German and
# Polish their own, real translations so the level-end script has proper
# localized text to display.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C31").Value = "Großartig!"
$ws.Range("D31").Value = "Świetnie!"

# Reflect the author's final scroll/selection state in the sheet view.
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("D32").Select()
